# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures on the active worksheet to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @{ D = new price text; E = new volume text }
# A "D" value of $null means the Price column is unchanged for that row.
$updates = @{
    2 = @{ D = "69.104.68"; E = "  -1.67%  " }
    3 = @{ D = "3.523.38"; E = "  -1.93%  " }
    4 = @{ D = $null; E = "  +0.13%  " }
    5 = @{ D = "570.71"; E = "  -1.36%  " }
    6 = @{ D = "181.57"; E = "  -4.87%  " }
    7 = @{ D = "3.518.74"; E = "  -1.85%  " }
    8 = @{ D = "0.613"; E = "  -3.24%  " }
    9 = @{ D = $null; E = "  +0.06%  " }
    10 = @{ D = $null; E = "  +6.12%  " }
    11 = @{ D = "0.636"; E = "  -3.89%  " }
    12 = @{ D = "53.50"; E = "  -5.57%  " }
    13 = @{ D = $null; E = "  +0.36%  " }
    14 = @{ D = "9.43"; E = "  -3.40%  " }
    15 = @{ D = "4.102.81"; E = "  -1.65%  " }
    16 = @{ D = "19.21"; E = "  -4.80%  " }
    17 = @{ D = "3.520.26"; E = "  -1.96%  " }
    18 = @{ D = "69.089.86"; E = "  -1.55%  " }
    19 = @{ D = "12.41"; E = "  -1.07%  " }
    20 = @{ D = $null; E = "  -1.39%  " }
    21 = @{ D = "538.31"; E = "  +13.45%  " }
    22 = @{ D = $null; E = "  -1.48%  " }
    23 = @{ D = "19.76"; E = "  +1.60%  " }
    24 = @{ D = "4.93"; E = "  -3.61%  " }
    25 = @{ D = "4.37"; E = "  -0.02%  " }
    26 = @{ D = "94.14"; E = "  +6.06%  " }
    27 = @{ D = "11.01"; E = "  -0.92%  " }
    28 = @{ D = "2.90"; E = "  -5.59%  " }
    29 = @{ D = "9.03"; E = "  -2.66%  " }
    30 = @{ D = "31.62"; E = "  -1.70%  " }
    31 = @{ D = "7.25"; E = "  -6.35%  " }
    32 = @{ D = "12.52"; E = "  +3.19%  " }
    33 = @{ D = "64.82"; E = "  -1.97%  " }
    34 = @{ D = "0.114"; E = "  -5.25%  " }
    35 = @{ D = "571.31"; E = "  -2.83%  " }
    36 = @{ D = "3.15"; E = "  +7.65%  " }
    37 = @{ D = "38.05"; E = "  -3.97%  " }
    38 = @{ D = $null; E = "  -0.17%  " }
    39 = @{ D = "0.397"; E = "  -0.75%  " }
    40 = @{ D = "0.0₃0761"; E = "  -5.77%  " }
    41 = @{ D = "3.36"; E = "  -5.59%  " }
    42 = @{ D = "0.132"; E = "  -7.17%  " }
    43 = @{ D = "3.05"; E = "  -2.98%  " }
    44 = @{ D = $null; E = "  +4.50%  " }
    45 = @{ D = "2.95"; E = "  -5.06%  " }
    46 = @{ D = "0.0440"; E = "  -1.68%  " }
    47 = @{ D = "3.160.92"; E = "  -2.18%  " }
    48 = @{ D = $null; E = "  -4.59%  " }
    49 = @{ D = "0.133"; E = "  -2.84%  " }
    50 = @{ D = "0.999"; E = "  -0.04%  " }
    51 = @{ D = "135.70"; E = "  -1.76%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($null -ne $vals.D) {
        # Force text storage so Excel does not reinterpret values such as
        # "69.104.68" or "0.999" as numbers, then restore the default style
        # so no residual formatting is left behind on the cell.
        $priceCell = $ws.Range("D" + $row)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $vals.D
        $priceCell.Style = "Normal"
    }

    $ws.Range("E" + $row).Value = $vals.E
}
